# Natmi following Dr Hou advice
# Update the LR-pair result table (Lrpap1-Vldlr) to the recomputed values,
# expanding the sending/target cluster grid from 3 clusters (12 combos) to
# the full 4-cluster grid (16 combos): ECs, FAPs, M2, sCs x ECs, FAPs, M2, sCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 16,20

$data[0,0] = "ECs"
$data[0,1] = "Lrpap1"
$data[0,2] = "Vldlr"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 8.191447666666667
$data[0,7] = 24.574343
$data[0,8] = 0.185794284429433
$data[0,9] = 0.185794284429433
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.463191
$data[0,13] = 1.389573
$data[0,14] = 0.0353316468093919
$data[0,15] = 0.0353316468093919
$data[0,16] = 3.794204836171
$data[0,17] = 34.147843525539
$data[0,18] = 0.006564418036664427
$data[0,19] = 0.006564418036664427

$data[1,0] = "ECs"
$data[1,1] = "Lrpap1"
$data[1,2] = "Vldlr"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 8.191447666666667
$data[1,7] = 24.574343
$data[1,8] = 0.185794284429433
$data[1,9] = 0.185794284429433
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 9.932929333333332
$data[1,13] = 29.798788
$data[1,14] = 0.7576717833204485
$data[1,15] = 0.7576717833204486
$data[1,16] = 81.36507081069821
$data[1,17] = 732.2856372962839
$data[1,18] = 0.1407710868143951
$data[1,19] = 0.1407710868143951

$data[2,0] = "ECs"
$data[2,1] = "Lrpap1"
$data[2,2] = "Vldlr"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 8.191447666666667
$data[2,7] = 24.574343
$data[2,8] = 0.185794284429433
$data[2,9] = 0.185794284429433
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.1285543333333333
$data[2,13] = 0.385663
$data[2,14] = 0.009805968382697785
$data[2,15] = 0.009805968382697785
$data[2,16] = 1.053046093823222
$data[2,17] = 9.477414844408999
$data[2,18] = 0.001821892878800979
$data[2,19] = 0.001821892878800979

$data[3,0] = "ECs"
$data[3,1] = "Lrpap1"
$data[3,2] = "Vldlr"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 8.191447666666667
$data[3,7] = 24.574343
$data[3,8] = 0.185794284429433
$data[3,9] = 0.185794284429433
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.585130333333333
$data[3,13] = 7.755391
$data[3,14] = 0.1971906014874617
$data[3,15] = 0.1971906014874618
$data[3,16] = 21.17595983701256
$data[3,17] = 190.583638533113
$data[3,18] = 0.03663688669957244
$data[3,19] = 0.03663688669957244

$data[4,0] = "FAPs"
$data[4,1] = "Lrpap1"
$data[4,2] = "Vldlr"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 17.317702
$data[4,7] = 51.95310600000001
$data[4,8] = 0.3927913821808575
$data[4,9] = 0.3927913821808576
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.463191
$data[4,13] = 1.389573
$data[4,14] = 0.0353316468093919
$data[4,15] = 0.0353316468093919
$data[4,16] = 8.021403707081999
$data[4,17] = 72.19263336373801
$data[4,18] = 0.01387796638498693
$data[4,19] = 0.01387796638498693

$data[5,0] = "FAPs"
$data[5,1] = "Lrpap1"
$data[5,2] = "Vldlr"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 17.317702
$data[5,7] = 51.95310600000001
$data[5,8] = 0.3927913821808575
$data[5,9] = 0.3927913821808576
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 9.932929333333332
$data[5,13] = 29.798788
$data[5,14] = 0.7576717833204485
$data[5,15] = 0.7576717833204486
$data[5,16] = 172.0155101817253
$data[5,17] = 1548.139591635528
$data[5,18] = 0.2976069470098742
$data[5,19] = 0.2976069470098742

$data[6,0] = "FAPs"
$data[6,1] = "Lrpap1"
$data[6,2] = "Vldlr"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 17.317702
$data[6,7] = 51.95310600000001
$data[6,8] = 0.3927913821808575
$data[6,9] = 0.3927913821808576
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.1285543333333333
$data[6,13] = 0.385663
$data[6,14] = 0.009805968382697785
$data[6,15] = 0.009805968382697785
$data[6,16] = 2.226265635475333
$data[6,17] = 20.036390719278
$data[6,18] = 0.003851699874661651
$data[6,19] = 0.003851699874661651

$data[7,0] = "FAPs"
$data[7,1] = "Lrpap1"
$data[7,2] = "Vldlr"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 17.317702
$data[7,7] = 51.95310600000001
$data[7,8] = 0.3927913821808575
$data[7,9] = 0.3927913821808576
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.585130333333333
$data[7,13] = 7.755391
$data[7,14] = 0.1971906014874617
$data[7,15] = 0.1971906014874618
$data[7,16] = 44.76851674382733
$data[7,17] = 402.916650694446
$data[7,18] = 0.07745476891133475
$data[7,19] = 0.07745476891133478

$data[8,0] = "M2"
$data[8,1] = "Lrpap1"
$data[8,2] = "Vldlr"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 12.35128266666667
$data[8,7] = 37.053848
$data[8,8] = 0.2801455637905346
$data[8,9] = 0.2801455637905346
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.463191
$data[8,13] = 1.389573
$data[8,14] = 0.0353316468093919
$data[8,15] = 0.0353316468093919
$data[8,16] = 5.721002969656
$data[8,17] = 51.489026726904
$data[8,18] = 0.009898004115065136
$data[8,19] = 0.009898004115065136

$data[9,0] = "M2"
$data[9,1] = "Lrpap1"
$data[9,2] = "Vldlr"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 12.35128266666667
$data[9,7] = 37.053848
$data[9,8] = 0.2801455637905346
$data[9,9] = 0.2801455637905346
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 9.932929333333332
$data[9,13] = 29.798788
$data[9,14] = 0.7576717833204485
$data[9,15] = 0.7576717833204486
$data[9,16] = 122.6844179040249
$data[9,17] = 1104.159761136224
$data[9,18] = 0.2122583889064868
$data[9,19] = 0.2122583889064868

$data[10,0] = "M2"
$data[10,1] = "Lrpap1"
$data[10,2] = "Vldlr"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 12.35128266666667
$data[10,7] = 37.053848
$data[10,8] = 0.2801455637905346
$data[10,9] = 0.2801455637905346
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.1285543333333333
$data[10,13] = 0.385663
$data[10,14] = 0.009805968382697785
$data[10,15] = 0.009805968382697785
$data[10,16] = 1.587810909024889
$data[10,17] = 14.290298181224
$data[10,18] = 0.002747098541083027
$data[10,19] = 0.002747098541083027

$data[11,0] = "M2"
$data[11,1] = "Lrpap1"
$data[11,2] = "Vldlr"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 12.35128266666667
$data[11,7] = 37.053848
$data[11,8] = 0.2801455637905346
$data[11,9] = 0.2801455637905346
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 2.585130333333333
$data[11,13] = 7.755391
$data[11,14] = 0.1971906014874617
$data[11,15] = 0.1971906014874618
$data[11,16] = 31.92967547717423
$data[11,17] = 287.367079294568
$data[11,18] = 0.05524207222789959
$data[11,19] = 0.0552420722278996

$data[12,0] = "sCs"
$data[12,1] = "Lrpap1"
$data[12,2] = "Vldlr"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 6.228371
$data[12,7] = 18.685113
$data[12,8] = 0.1412687695991749
$data[12,9] = 0.1412687695991749
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 0.463191
$data[12,13] = 1.389573
$data[12,14] = 0.0353316468093919
$data[12,15] = 0.0353316468093919
$data[12,16] = 2.884925391861
$data[12,17] = 25.964328526749
$data[12,18] = 0.004991258272675407
$data[12,19] = 0.004991258272675407

$data[13,0] = "sCs"
$data[13,1] = "Lrpap1"
$data[13,2] = "Vldlr"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 6.228371
$data[13,7] = 18.685113
$data[13,8] = 0.1412687695991749
$data[13,9] = 0.1412687695991749
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 9.932929333333332
$data[13,13] = 29.798788
$data[13,14] = 0.7576717833204485
$data[13,15] = 0.7576717833204486
$data[13,16] = 61.86596900478266
$data[13,17] = 556.793721043044
$data[13,18] = 0.1070353605896924
$data[13,19] = 0.1070353605896924

$data[14,0] = "sCs"
$data[14,1] = "Lrpap1"
$data[14,2] = "Vldlr"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 6.228371
$data[14,7] = 18.685113
$data[14,8] = 0.1412687695991749
$data[14,9] = 0.1412687695991749
$data[14,10] = 2
$data[14,11] = 0.6666666666666666
$data[14,12] = 0.1285543333333333
$data[14,13] = 0.385663
$data[14,14] = 0.009805968382697785
$data[14,15] = 0.009805968382697785
$data[14,16] = 0.8006840816576666
$data[14,17] = 7.206156734919
$data[14,18] = 0.001385277088152127
$data[14,19] = 0.001385277088152127

$data[15,0] = "sCs"
$data[15,1] = "Lrpap1"
$data[15,2] = "Vldlr"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 6.228371
$data[15,7] = 18.685113
$data[15,8] = 0.1412687695991749
$data[15,9] = 0.1412687695991749
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 2.585130333333333
$data[15,13] = 7.755391
$data[15,14] = 0.1971906014874617
$data[15,15] = 0.1971906014874618
$data[15,16] = 16.10115079935367
$data[15,17] = 144.910357194183
$data[15,18] = 0.02785687364865494
$data[15,19] = 0.02785687364865494

$ws.Range("A2:T17").Value = $data
